$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.457.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +9.33%  '

$ws.Range("D3").Value = "'3.466.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.47%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'414.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.19%  '

$ws.Range("D6").Value = "'123.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.94%  '

$ws.Range("D7").Value = "'3.461.95"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.45%  '

$ws.Range("D8").Value = "'0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.62%  '

$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").Value = "'0.665"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.34%  '

$ws.Range("D11").Value = "'0.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +33.92%  '

$ws.Range("D12").Value = "'41.36"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").Value = "'3.999.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.99%  '

$ws.Range("D15").Value = "'8.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.33%  '

$ws.Range("D16").Value = "'19.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.31%  '

$ws.Range("D17").Value = "'3.459.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.41%  '

$ws.Range("D18").Value = "'62.289.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.43%  '

$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").Value = "'11.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("E21").Value = '  +22.63%  '

$ws.Range("E22").Value = '  +1.30%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = "'318.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.39%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = "'82.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.59%  '

$ws.Range("E25").Value = '  +1.07%  '

$ws.Range("D26").Value = "'3.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").Value = "'31.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.50%  '

$ws.Range("D28").Value = "'7.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.53%  '

$ws.Range("D29").Value = "'7.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").Value = "'4.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("E31").Value = '  +2.53%  '

$ws.Range("E32").Value = '  +3.55%  '

$ws.Range("D33").Value = "'11.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.41%  '

$ws.Range("D34").Value = "'42.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.21%  '

$ws.Range("D35").Value = "'2.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +19.92%  '

$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("D37").Value = "'0.0486"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").Value = "'52.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("E39").Value = '  +1.15%  '

$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").Value = "'3.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.15%  '

$ws.Range("E42").Value = '  +8.43%  '

$ws.Range("D43").Value = "'0.126"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.44%  '

$ws.Range("D44").Value = "'134.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").Value = "'17.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.09%  '

$ws.Range("D46").Value = "'0.285"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.39%  '

$ws.Range("D47").Value = "'3.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.29%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'2.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.87%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'22.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.40%  '

$ws.Range("D50").Value = "'2.210.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.12%  '

$ws.Range("D51").Value = "'3.800.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.06%  '
